$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '307.24'
$r.ClearFormats()
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '0.93%'
$r.ClearFormats()
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '36.30'
$r.ClearFormats()
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '1.18%'
$r.ClearFormats()
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '5.055'
$r.ClearFormats()
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '-0.39%'
$r.ClearFormats()
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '0.08079'
$r.ClearFormats()
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '0.34%'
$r.ClearFormats()
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '2.160'
$r.ClearFormats()
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '12.27%'
$r.ClearFormats()
$r = $ws.Range("B7")
$r.NumberFormat = "@"
$r.Value = 'GateToken'
$r.ClearFormats()
$r = $ws.Range("C7")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$r.ClearFormats()
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '4.147'
$r.ClearFormats()
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '-0.05%'
$r.ClearFormats()
$r = $ws.Range("B8")
$r.NumberFormat = "@"
$r.Value = 'KuCoinToken'
$r.ClearFormats()
$r = $ws.Range("C8")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$r.ClearFormats()
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '7.840'
$r.ClearFormats()
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '-0.12%'
$r.ClearFormats()
$r = $ws.Range("B9")
$r.NumberFormat = "@"
$r.Value = 'MXToken'
$r.ClearFormats()
$r = $ws.Range("C9")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r.ClearFormats()
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.9276'
$r.ClearFormats()
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '-0.38%'
$r.ClearFormats()
$r = $ws.Range("B10")
$r.NumberFormat = "@"
$r.Value = 'LiechtensteinCryptoassetsExchange'
$r.ClearFormats()
$r = $ws.Range("C10")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$r.ClearFormats()
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.1427'
$r.ClearFormats()
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '11.73%'
$r.ClearFormats()
$r = $ws.Range("B11")
$r.NumberFormat = "@"
$r.Value = 'WazirX'
$r.ClearFormats()
$r = $ws.Range("C11")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$r.ClearFormats()
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.1929'
$r.ClearFormats()
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '0.71%'
$r.ClearFormats()
$r = $ws.Range("B12")
$r.NumberFormat = "@"
$r.Value = 'MandalaExchangeToken'
$r.ClearFormats()
$r = $ws.Range("C12")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$r.ClearFormats()
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.09099'
$r.ClearFormats()
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '-0.43%'
$r.ClearFormats()
$r = $ws.Range("B13")
$r.NumberFormat = "@"
$r.Value = 'BitrueCoin'
$r.ClearFormats()
$r = $ws.Range("C13")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$r.ClearFormats()
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.03454'
$r.ClearFormats()
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '-0.93%'
$r.ClearFormats()
$r = $ws.Range("B14")
$r.NumberFormat = "@"
$r.Value = 'BitMartToken'
$r.ClearFormats()
$r = $ws.Range("C14")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$r.ClearFormats()
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.09912'
$r.ClearFormats()
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '0.08%'
$r.ClearFormats()
$r = $ws.Range("B15")
$r.NumberFormat = "@"
$r.Value = 'BitForexToken'
$r.ClearFormats()
$r = $ws.Range("C15")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$r.ClearFormats()
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.001405'
$r.ClearFormats()
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '-1.53%'
$r.ClearFormats()
$r = $ws.Range("B16")
$r.NumberFormat = "@"
$r.Value = 'TigerCash'
$r.ClearFormats()
$r = $ws.Range("C16")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$r.ClearFormats()
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '0.006316'
$r.ClearFormats()
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '-4.96%'
$r.ClearFormats()
$r = $ws.Range("B17")
$r.NumberFormat = "@"
$r.Value = 'LEO'
$r.ClearFormats()
$r = $ws.Range("C17")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$r.ClearFormats()
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '3.842'
$r.ClearFormats()
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '6.34%'
$r.ClearFormats()
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '3.398'
$r.ClearFormats()
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '7.23%'
$r.ClearFormats()
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.3454'
$r.ClearFormats()
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '0.97%'
$r.ClearFormats()
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '0.1302'
$r.ClearFormats()
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '-2.58%'
$r.ClearFormats()
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '4.835'
$r.ClearFormats()
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '-6.66%'
$r.ClearFormats()
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '-7.68%'
$r.ClearFormats()
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '0.04361'
$r.ClearFormats()
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '-0.98%'
$r.ClearFormats()
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '0.001231'
$r.ClearFormats()
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '-0.48%'
$r.ClearFormats()
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.004303'
$r.ClearFormats()
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '0.0001299'
$r.ClearFormats()
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '-0.31%'
$r.ClearFormats()
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.02020'
$r.ClearFormats()
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '1.68%'
$r.ClearFormats()
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.05165'
$r.ClearFormats()
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '0.32%'
$r.ClearFormats()
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.007522'
$r.ClearFormats()
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '-1.01%'
$r.ClearFormats()
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.01016'
$r.ClearFormats()
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '0.99%'
$r.ClearFormats()
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.1366'
$r.ClearFormats()
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '0.12%'
$r.ClearFormats()
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.002148'
$r.ClearFormats()
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '2.06%'
$r.ClearFormats()
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.009963'
$r.ClearFormats()
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '-6.93%'
$r.ClearFormats()
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.00006287'
$r.ClearFormats()
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '-3.11%'
$r.ClearFormats()
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '-0.29%'
$r.ClearFormats()
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '64.85'
$r.ClearFormats()
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '-0.16%'
$r.ClearFormats()
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.001250'
$r.ClearFormats()
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '-22.03%'
$r.ClearFormats()
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.00002098'
$r.ClearFormats()
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '-0.29%'
$r.ClearFormats()
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.0001998'
$r.ClearFormats()
